# Updated capital structure database
# Apply updated values to rows 2 and 3 for the Nigeria broadcasting dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = -0.133
    "G2" = -0.1591304347826087
    "H2" = -0.1591304347826087
    "I2" = -0.3695652173913043
    "J2" = -0.3695652173913043
    "K2" = -3.72
    "L2" = -0.3234782608695652
    "U2" = 0.273
    "V2" = 0.04340222575516693
    "W2" = -0.1917525773195876
    "X2" = 0.100397840228604
    "Y2" = -0.2921504175481917
    "Z2" = 0.5155563525508833
    "AA2" = -0.1905316955079351
    "AB2" = 0.08611495961164105
    "AC2" = -0.2766466551195761
    "AD2" = 2.89
    "AF2" = 2.89
    "AG2" = 2.617
    "AH2" = 0.3148148148148148
    "AI2" = 0.1731575793888556
    "AJ2" = 0.293813854271921
    "AK2" = 0.1594079308034355
    "AL2" = 0.294
    "AM2" = 0.294
    "AN2" = -1.090566037735849
    "AO2" = -14.45578231292517
    "AP2" = -0.9875471698113207
    "AQ2" = -14.45578231292517

    "D3" = -0.133
    "G3" = -0.1591304347826087
    "H3" = -0.1591304347826087
    "I3" = -0.3695652173913043
    "J3" = -0.3695652173913043
    "K3" = -3.72
    "L3" = -0.3234782608695652
    "U3" = 0.273
    "V3" = 0.04340222575516693
    "W3" = -0.1917525773195876
    "X3" = 0.100397840228604
    "Y3" = -0.2921504175481917
    "Z3" = 0.5155563525508833
    "AA3" = -0.1905316955079351
    "AB3" = 0.08611495961164105
    "AC3" = -0.2766466551195761
    "AD3" = 2.89
    "AF3" = 2.89
    "AG3" = 2.617
    "AH3" = 0.3148148148148148
    "AI3" = 0.1731575793888556
    "AJ3" = 0.293813854271921
    "AK3" = 0.1594079308034355
    "AL3" = 0.294
    "AM3" = 0.294
    "AN3" = -1.090566037735849
    "AO3" = -14.45578231292517
    "AP3" = -0.9875471698113207
    "AQ3" = -14.45578231292517
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
